$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.300385713577271
$ws.Range("B1").Value = 1.93176805973053
$ws.Range("C1").Value = 5.213043212890625
$ws.Range("D1").Value = 1.979284882545471
$ws.Range("E1").Value = 1.089383721351624
